# Task: Completed daily operations, 8 hours, 09/21
# Add a new time-log entry row (row 6) mirroring the existing rows:
#   A: Date, B: Name of Task ("Internship"), C: Description (completed task text)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry for 09/21/2023 (serial date 45190, matching the other date cells)
$ws.Range("A6").Value = 45190
$ws.Range("B6").Value = "Internship"
$ws.Range("C6").Value = "Completed 8 hours assisting with daily operations"

# Match the existing date-column formatting (style used by A2:A5)
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat

# Move the active selection to D6, matching the post-edit workbook state
$ws.Range("D6").Select()
